# Updated cryptos list - applies Price (D) and Volume(1h) (E) changes,
# plus the WEMIXToken/LidoDAOToken row swap (rows 33-34 incl. Coin/Link).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.983.58'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '2.575.74'
$ws.Range("E3").Value = '  +2.63%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'" + '302.93'
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").Value = "'" + '97.38'
$ws.Range("E6").Value = '  +4.59%  '
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").Value = "'" + '0.550'
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("D10").Value = "'" + '36.55'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("E11").Value = '  +1.29%  '
$ws.Range("D12").Value = "'" + '7.76'
$ws.Range("E12").Value = '  +1.55%  '
$ws.Range("E13").Value = '  +6.90%  '
$ws.Range("D14").Value = '2.540.54'
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("D15").Value = "'" + '0.889'
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").Value = "'" + '14.37'
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("D17").Value = '43.020.42'
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("D18").Value = "'" + '13.00'
$ws.Range("E18").Value = '  +6.47%  '
$ws.Range("D19").Value = '0.0₃0997'
$ws.Range("E19").Value = '  +4.15%  '
$ws.Range("D20").Value = "'" + '6.64'
$ws.Range("E20").Value = '  +2.44%  '
$ws.Range("D21").Value = "'" + '72.00'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").Value = "'" + '255.16'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").Value = "'" + '2.98'
$ws.Range("E23").Value = '  +3.62%  '
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").Value = "'" + '28.70'
$ws.Range("E25").Value = '  +1.61%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = "'" + '10.28'
$ws.Range("E27").Value = '  +3.10%  '
$ws.Range("D28").Value = "'" + '37.85'
$ws.Range("E28").Value = '  +3.25%  '
$ws.Range("E29").Value = '  -4.02%  '
$ws.Range("D30").Value = "'" + '6.06'
$ws.Range("E30").Value = '  +1.79%  '
$ws.Range("D31").Value = "'" + '155.57'
$ws.Range("E31").Value = '  +3.43%  '
$ws.Range("E32").Value = '  +1.88%  '
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = "'" + '2.76'
$ws.Range("E33").Value = '  +1.38%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = "'" + '3.39'
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("D35").Value = "'" + '0.0814'
$ws.Range("E35").Value = '  +2.89%  '
$ws.Range("D36").Value = "'" + '18.46'
$ws.Range("E36").Value = '  +11.51%  '
$ws.Range("E37").Value = '  +1.24%  '
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("D39").Value = "'" + '23.81'
$ws.Range("E39").Value = '  +1.15%  '
$ws.Range("D40").Value = "'" + '3.43'
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("E41").Value = '  +1.36%  '
$ws.Range("D42").Value = "'" + '3.89'
$ws.Range("E42").Value = '  +2.85%  '
$ws.Range("E43").Value = '  +24.39%  '
$ws.Range("D44").Value = '2.068.65'
$ws.Range("E44").Value = '  +2.89%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("E46").Value = '  +4.64%  '
$ws.Range("D47").Value = "'" + '85.44'
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("D48").Value = "'" + '77.26'
$ws.Range("E48").Value = '  +15.08%  '
$ws.Range("D49").Value = "'" + '106.34'
$ws.Range("E49").Value = '  +3.81%  '
$ws.Range("D50").Value = '2.822.96'
$ws.Range("E50").Value = '  +2.26%  '
$ws.Range("E51").Value = '  +2.96%  '
